# Weekly refresh of Fruta/Hortaliza prices for
# "Vega Central Mapocho de Santiago - Papaya": the per-row Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio ponderado, Unidad de
# comercializacion, Precio $/Kg and Kg/unidad values are rotated across
# rows 2-15 as the underlying daily records roll forward to the next
# reporting week. Mercado/Region/Codreg/Tipo/Producto/Categoria/Variedad/
# Origen columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=44309; L="Primera"; M=10;  N=1600;  O=1600;  P=1600;  Q="`$/kilo (en caja de 15 kilos)"; S=1600; T=1 },
    @{ Row=3;  D=44195; L="Primera"; M=20;  N=15000; O=15000; P=15000; Q="`$/bandeja 10 kilos";           S=1500; T=10 },
    @{ Row=4;  D=44336; L="Primera"; M=10;  N=1500;  O=1500;  P=1500;  Q="`$/kilo (en caja de 15 kilos)"; S=1500; T=1 },
    @{ Row=5;  D=44391; L="Primera"; M=15;  N=1500;  O=1500;  P=1500;  Q="`$/kilo (en caja de 15 kilos)"; S=1500; T=1 },
    @{ Row=6;  D=44391; L="Segunda"; M=20;  N=1000;  O=1000;  P=1000;  Q="`$/kilo (en caja de 15 kilos)"; S=1000; T=1 },
    @{ Row=7;  D=44371; L="Primera"; M=20;  N=1800;  O=1800;  P=1800;  Q="`$/kilo (en caja de 15 kilos)"; S=1800; T=1 },
    @{ Row=8;  D=44371; L="Segunda"; M=30;  N=1200;  O=1200;  P=1200;  Q="`$/kilo (en caja de 15 kilos)"; S=1200; T=1 },
    @{ Row=9;  D=44400; L="Primera"; M=25;  N=1500;  O=1500;  P=1500;  Q="`$/kilo (en caja de 15 kilos)"; S=1500; T=1 },
    @{ Row=10; D=44343; L="Primera"; M=20;  N=1700;  O=1700;  P=1700;  Q="`$/kilo (en caja de 15 kilos)"; S=1700; T=1 },
    @{ Row=11; D=44904; L="Primera"; M=45;  N=15000; O=15000; P=15000; Q="`$/bandeja 10 kilos";           S=1500; T=10 },
    @{ Row=12; D=44904; L="Segunda"; M=60;  N=10000; O=10000; P=10000; Q="`$/bandeja 10 kilos";           S=1000; T=10 },
    @{ Row=13; D=44292; L="Primera"; M=50;  N=14000; O=14000; P=14000; Q="`$/bandeja 10 kilos";           S=1400; T=10 },
    @{ Row=14; D=44880; L="Primera"; M=200; N=20000; O=20000; P=20000; Q="`$/bandeja 10 kilos";           S=2000; T=10 },
    @{ Row=15; D=44880; L="Segunda"; M=180; N=15000; O=15000; P=15000; Q="`$/bandeja 10 kilos";           S=1500; T=10 }
)

foreach ($r in $rows) {
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("L$($r.Row)").Value = $r.L
    $ws.Range("M$($r.Row)").Value = $r.M
    $ws.Range("N$($r.Row)").Value = $r.N
    $ws.Range("O$($r.Row)").Value = $r.O
    $ws.Range("P$($r.Row)").Value = $r.P
    $ws.Range("Q$($r.Row)").Value = $r.Q
    $ws.Range("S$($r.Row)").Value = $r.S
    $ws.Range("T$($r.Row)").Value = $r.T
}
